# Update the "想去人数" (people interested) counts in column F
# for the "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 3
$ws1.Range("F7").Value  = 1778
$ws1.Range("F11").Value = 2084
$ws1.Range("F12").Value = 27
$ws1.Range("F13").Value = 144
$ws1.Range("F14").Value = 1338
$ws1.Range("F15").Value = 470
$ws1.Range("F23").Value = 56
$ws1.Range("F24").Value = 13
$ws1.Range("F25").Value = 1123
$ws1.Range("F27").Value = 336
$ws1.Range("F29").Value = 272

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 3
$ws4.Range("F7").Value  = 1778
$ws4.Range("F12").Value = 2084
$ws4.Range("F13").Value = 27
$ws4.Range("F14").Value = 144
$ws4.Range("F15").Value = 1338
$ws4.Range("F16").Value = 470
$ws4.Range("F24").Value = 56
$ws4.Range("F25").Value = 13
$ws4.Range("F26").Value = 1123
$ws4.Range("F28").Value = 336
$ws4.Range("F30").Value = 272

$wb.Save()
